# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled
# update). Only the Price (D) and Volume(1h) (E) columns move for most
# rows; rows 12/13 (Polkadot / WrappedEther) also swap position in the
# ranking, so their Coin name (B) and Link (C) are updated too.
#
# Price values are written with NumberFormat "@" wherever the new text
# looks like a plain number (e.g. "3.200", "1.001") so Excel keeps
# storing them as text - exactly like the source data - instead of
# silently coercing them into numbers and dropping significant trailing
# zeros (e.g. "3.200" -> 3.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.919.39'
$ws.Range("E2").Value = '  -0.69%  '

# Row 3
$ws.Range("D3").Value = '1.639.36'
$ws.Range("E3").Value = '  -0.65%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.42%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.12'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5054'
$ws.Range("E6").Value = '  -0.89%  '

# Row 7
$ws.Range("E7").Value = '  -1.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06452'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2574'
$ws.Range("E9").Value = '  -0.53%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.54'
$ws.Range("E10").Value = '  -0.44%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07739'
$ws.Range("E11").Value = '  -0.49%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.255'
$ws.Range("E12").Value = '  -0.29%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.635.01'
$ws.Range("E13").Value = '  -0.96%  '

# Row 14
$ws.Range("D14").Value = '1.864.99'
$ws.Range("E14").Value = '  -0.74%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5464'
$ws.Range("E15").Value = '  -0.03%  '

# Row 16
$ws.Range("D16").Value = '0.0₅7914'
$ws.Range("E16").Value = '  -0.68%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.26'
$ws.Range("E17").Value = '  +0.90%  '

# Row 18
$ws.Range("D18").Value = '25.919.22'
$ws.Range("E18").Value = '  -0.73%  '

# Row 19
$ws.Range("E19").Value = '  -1.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.34'
$ws.Range("E20").Value = '  -2.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.375'
$ws.Range("E21").Value = '  -0.42%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.885'
$ws.Range("E22").Value = '  -1.63%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.973'
$ws.Range("E23").Value = '  -0.78%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -1.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.862'
$ws.Range("E25").Value = '  +0.36%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.79'
$ws.Range("E26").Value = '  -2.26%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1136'
$ws.Range("E27").Value = '  -2.68%  '

# Row 28
$ws.Range("E28").Value = '  -0.75%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.778'
$ws.Range("E29").Value = '  -2.10%  '

# Row 30
$ws.Range("E30").Value = '  -0.09%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04921'
$ws.Range("E31").Value = '  -3.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.268'
$ws.Range("E32").Value = '  -1.90%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.200'
$ws.Range("E33").Value = '  -0.99%  '

# Row 34
$ws.Range("E34").Value = '  +0.04%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.362'
$ws.Range("E35").Value = '  -0.22%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8916'
$ws.Range("E36").Value = '  -2.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.622'
$ws.Range("E37").Value = '  -3.17%  '

# Row 38
$ws.Range("D38").Value = '1.149.33'
$ws.Range("E38").Value = '  -1.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5589'
$ws.Range("E39").Value = '  -1.57%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01562'
$ws.Range("E40").Value = '  -0.82%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9995'
$ws.Range("E41").Value = '  -1.40%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.694'
$ws.Range("E42").Value = '  +0.50%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.65'
$ws.Range("E43").Value = '  -0.48%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8042'
$ws.Range("E44").Value = '  -2.45%  '

# Row 45
$ws.Range("D45").Value = '1.776.20'
$ws.Range("E45").Value = '  -0.83%  '

# Row 46
$ws.Range("D46").Value = '0.0₈116'
$ws.Range("E46").Value = '  +3.73%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4515'
$ws.Range("E47").Value = '  -1.03%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  -0.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.63'
$ws.Range("E49").Value = '  -1.08%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05047'
$ws.Range("E50").Value = '  -0.79%  '

# Row 51
$ws.Range("E51").Value = '  -1.02%  '
